$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Canada Premier League - Valour vs Vancouver FC) - updated odds
$ws.Range("Y3").Value = 32
$ws.Range("AA3").Value = 7.1
$ws.Range("AI3").Value = 19

# Row 6 (Colombia Primera B - Patriotas vs Real Cundinamarca) - updated odds
$ws.Range("G6").Value = 1.65
$ws.Range("H6").Value = 3.65
$ws.Range("I6").Value = 4.85
$ws.Range("L6").Value = 1.27
$ws.Range("M6").Value = 3.1
$ws.Range("N6").Value = 1.8
$ws.Range("O6").Value = 1.8
$ws.Range("P6").Value = 1.38
$ws.Range("Q6").Value = 2.6
$ws.Range("R6").Value = 1.8
$ws.Range("S6").Value = 1.8
$ws.Range("T6").Value = 6.6
$ws.Range("U6").Value = 7.5
$ws.Range("X6").Value = 13.5
$ws.Range("Y6").Value = 27
$ws.Range("Z6").Value = 10.25
$ws.Range("AA6").Value = 7.2
$ws.Range("AB6").Value = 16.5
$ws.Range("AC6").Value = 80
$ws.Range("AG6").Value = 15.5

# Row 7 (Ecuador Serie B - Gualaceo vs 22 de Julio) - updated odds (AD7 gains a value, was empty)
$ws.Range("H7").Value = 4.4
$ws.Range("I7").Value = 7.6
$ws.Range("L7").Value = 1.25
$ws.Range("M7").Value = 3.25
$ws.Range("N7").Value = 1.75
$ws.Range("O7").Value = 1.87
$ws.Range("R7").Value = 2.02
$ws.Range("S7").Value = 1.62
$ws.Range("T7").Value = 6.3
$ws.Range("U7").Value = 6.1
$ws.Range("V7").Value = 8.5
$ws.Range("W7").Value = 8.5
$ws.Range("Z7").Value = 10.75
$ws.Range("AA7").Value = 8.75
$ws.Range("AD7").Value = 700
$ws.Range("AE7").Value = 17.5
$ws.Range("AF7").Value = 50
$ws.Range("AG7").Value = 25
$ws.Range("AH7").Value = 200
$ws.Range("AI7").Value = 100
$ws.Range("AJ7").Value = 90

# Row 8 (Ecuador Serie B - Vargas Torres vs Atletico Vinotinto) - updated odds
$ws.Range("H8").Value = 3.15
$ws.Range("I8").Value = 2.52
$ws.Range("L8").Value = 1.32
$ws.Range("M8").Value = 2.87
$ws.Range("N8").Value = 1.93
$ws.Range("O8").Value = 1.7
$ws.Range("P8").Value = 1.39
$ws.Range("Q8").Value = 2.57
$ws.Range("R8").Value = 1.7
$ws.Range("S8").Value = 1.9
$ws.Range("T8").Value = 8.25
$ws.Range("U8").Value = 13.5
$ws.Range("X8").Value = 23
$ws.Range("Y8").Value = 32
$ws.Range("Z8").Value = 9
$ws.Range("AA8").Value = 6.1
$ws.Range("AB8").Value = 13.5
$ws.Range("AC8").Value = 65
$ws.Range("AD8").Value = 500
$ws.Range("AE8").Value = 8
$ws.Range("AF8").Value = 12.5
$ws.Range("AG8").Value = 9.5
$ws.Range("AI8").Value = 21
$ws.Range("AJ8").Value = 30

# Row 13 (Uruguay Segunda Division - Central Esp. vs Colon) - updated odds
$ws.Range("G13").Value = 3.85
$ws.Range("H13").Value = 3.25
$ws.Range("I13").Value = 1.9
$ws.Range("L13").Value = 1.42
$ws.Range("M13").Value = 2.45
$ws.Range("R13").Value = 2.05
$ws.Range("S13").Value = 1.6
$ws.Range("T13").Value = 8.75
$ws.Range("U13").Value = 19
$ws.Range("V13").Value = 14
$ws.Range("W13").Value = 60
$ws.Range("X13").Value = 45
$ws.Range("Y13").Value = 60
$ws.Range("Z13").Value = 7.3
$ws.Range("AA13").Value = 6.5
$ws.Range("AB13").Value = 19.5
$ws.Range("AE13").Value = 5.6
$ws.Range("AF13").Value = 7.8
$ws.Range("AG13").Value = 9
$ws.Range("AH13").Value = 15.5
$ws.Range("AI13").Value = 18

# Row 14 (Uruguay Segunda Division - Atenas vs Maldonado) - updated odds
$ws.Range("G14").Value = 3.6
$ws.Range("H14").Value = 3.15
$ws.Range("I14").Value = 2.02
$ws.Range("M14").Value = 2.35
$ws.Range("N14").Value = 2.35
$ws.Range("O14").Value = 1.47
$ws.Range("T14").Value = 8
$ws.Range("U14").Value = 17
$ws.Range("V14").Value = 13.5
$ws.Range("W14").Value = 55
$ws.Range("Z14").Value = 6.9
$ws.Range("AA14").Value = 6.3
$ws.Range("AE14").Value = 5.5
$ws.Range("AF14").Value = 8.25
$ws.Range("AG14").Value = 9.25
$ws.Range("AH14").Value = 17.5
$ws.Range("AI14").Value = 20

# Row 16 (USA USL Championship - Birmingham vs Indy Eleven) - updated odds
$ws.Range("H16").Value = 3.5
$ws.Range("T16").Value = 10.5
$ws.Range("U16").Value = 14.5
$ws.Range("AE16").Value = 10.5
